# Generate Report for Handoff
# Replace the old localization run's identifying GUID / hashes / timestamps
# with the new ones, across all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "df7c96cc-347c-4db0-9f5b-a1fb97e5ba8d"
$newGuid = "c4715556-fd27-4ab9-9547-a4e887f1fabf"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc62b239603f3c77e3f854f40cf7834acd574bbf/e2e/"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet 1): A2 file name, B2 path+name (hyperlink), G2 date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = ($newGuid + ".md")
$wsOverview.Range("B2").Value2 = ("e2e\" + $newGuid + ".md")
$wsOverview.Range("G2").Value2 = "2016-08-22 03:07:12"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($repoBase + $newGuid + ".md"), "", "", ("e2e\" + $newGuid + ".md"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet 2): A2 file name (hyperlink), G2 handoff file, H2 date
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = ($newGuid + ".md")
$wsZhCn.Range("G2").Value2 = ($newGuid + ".afebbfc1b9a11c84dfd2e988318118a693ae99a5.zh-cn.xlf")
$wsZhCn.Range("H2").Value2 = "2016-08-22 03:07:07"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($repoBase + $newGuid + ".md"), "", "", ($newGuid + ".md"))

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet 3): A2 file name (hyperlink), G2 handback file
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = ($newGuid + ".md")
$wsDeDe.Range("G2").Value2 = ($newGuid + ".afebbfc1b9a11c84dfd2e988318118a693ae99a5.de-de.xlf")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($repoBase + $newGuid + ".md"), "", "", ($newGuid + ".md"))
